# "Se añadio a Zulma" -- add a new response row for Zulma Clara Rios Rocha
# (grupo 3) to the bottom of the "Respuestas de formulario 1" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Respuestas de formulario 1")

$lastRow = 23
$newRow  = $lastRow + 1

# Bring the formatting (borders/font) of the last existing data row down to
# the new row before filling in the values, so the new row looks the same
# as the rest of the table.
$ws.Range("A" + $lastRow + ":C" + $lastRow).Copy($ws.Range("A" + $newRow + ":C" + $newRow))
$excel.CutCopyMode = 0
$ws.Rows($newRow).RowHeight = $ws.Rows($lastRow).RowHeight

# Fill in the new response: Paterno, Materno, Nombres, grupo
$ws.Range("A" + $newRow).Value = "Rios"
$ws.Range("B" + $newRow).Value = "Rocha"
$ws.Range("C" + $newRow).Value = "Zulma Clara"
$ws.Range("D" + $newRow).Value = 3

# Leave the selection where the author left it after adding the row.
$ws.Range("B27").Select()
